$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "telecomm/FFT/runme_large.sh"
$ws.Range("B18").Value = 0.07
$ws.Range("C18").Value = 0.05
$ws.Range("D18").Value = 0

$ws.Range("A19").Value = "telecomm/adpcm/runme_large.sh"
$ws.Range("B19").Value = 0.22
$ws.Range("C19").Value = 0.14
$ws.Range("D19").Value = 0.7

$ws.Range("F35").Select()
